$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalDistrbutionPayment")

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Fund *"
$ws.Range("B1").Value = "Investor *"
$ws.Range("C1").Value = "Capital Distribution *"
$ws.Range("D1").Value = "Amount *"
$ws.Range("E1").Value = "Cost Of Investment *"
$ws.Range("F1").Value = "Payment Date *"
$ws.Range("G1").Value = "Completed"
$ws.Range("H1").Value = "Folio No"

# --- Row 2 ---
$ws.Range("A2").Value = "SAAS Fund"
$ws.Range("B2").Value = "Investor 1"
$ws.Range("C2").Value = "Distribution 1 "
$ws.Range("G2").Value = "No"
$ws.Range("H2").Value = 6

# --- Row 3 ---
$ws.Range("A3").Value = "SAAS Fund"
$ws.Range("B3").Value = "Investor 2"
$ws.Range("C3").Value = "Distribution 1     "
$ws.Range("G3").Value = "No"
$ws.Range("H3").Value = 7

# --- Row 4 ---
$ws.Range("A4").Value = "SAAS Fund"
$ws.Range("B4").Value = "Investor 3"
$ws.Range("C4").Value = "Distribution 1"
$ws.Range("G4").Value = "No"
$ws.Range("H4").Value = 8

# --- Row 5 ---
$ws.Range("A5").Value = "SAAS Fund"
$ws.Range("B5").Value = "Investor 4"
$ws.Range("C5").Value = "Distribution 1"
$ws.Range("G5").Value = "Yes"
$ws.Range("H5").Value = 9

# --- Row 6 ---
$ws.Range("A6").Value = "SAAS Fund"
$ws.Range("B6").Value = "Investor 1"
$ws.Range("C6").Value = "Distribution 2"
$ws.Range("G6").Value = "Yes"
$ws.Range("H6").Value = 6

# --- Row 7 ---
$ws.Range("A7").Value = "SAAS Fund"
$ws.Range("B7").Value = "Investor 2"
$ws.Range("C7").Value = "Distribution 2"
$ws.Range("G7").Value = "Yes"
$ws.Range("H7").Value = 7

# Update selection to match the post-edit cursor position
$ws.Range("C4").Select()
